$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("Bough rice cooker" purchase): fold in the returned-item credit (row 28)
# directly into this row's cost as a formula, and update the comment to reflect
# that the item was purchased then later returned.
$ws.Range("D7").Formula = "=-287.52+68.23"
$ws.Range("F7").Value = "Bough rice cooker, then returned later."

# Row 28 was the separate "Returned the rice cooker" / Cookware credit entry.
# Its value has now been merged into row 7, so remove the now-redundant row,
# which shifts every following row up by one.
$ws.Rows(28).Delete()

# Reflect the cell the author was last looking at after the cleanup.
$ws.Range("F8").Select()
